$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) for the new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, border, center) from an existing header cell (AC1) to the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill data rows 2-48 with the team record values
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 83   # AD
    $ws.Cells.Item($r, 31).Value = 79   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
